# Remove the now-superseded "【会员购严选】苏州·二次元开放式年会- I COME ACG"
# listing (old row 5) from both the "展览" and "全部类型" sheets (they carry
# duplicate data), shifting every row below it up by one, and refresh the
# "想去人数" (want-to-go count, column F) figures that moved since the last
# scrape. Column A keeps its original per-row sequence number untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Destination rows 5..14 hold text-like values (plain "YYYY.MM.DD" dates
    # in column B, numeric-looking price strings in column G) that Excel
    # would otherwise auto-convert to dates/numbers on assignment. Force
    # those two columns to Text format first so the shifted values round
    # trip as the original strings.
    $ws.Range("B5:B14").NumberFormat = "@"
    $ws.Range("G5:G14").NumberFormat = "@"

    # Shift columns B:I of rows 6-15 up into rows 5-14 (this removes row 5's
    # old content and closes the gap), leaving column A's literal numbering
    # untouched.
    for ($r = 6; $r -le 15; $r++) {
        $destRow = $r - 1
        for ($c = 2; $c -le 9; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $ws.Cells.Item($r, $c).Value()
        }
    }

    # The old last row (15) is now a vacated duplicate; clear it entirely so
    # the sheet's used range shrinks back to row 14.
    $ws.Range("A15:I15").Clear()

    # Refresh the "想去人数" (F column) figures that changed between scrapes.
    $ws.Range("F2").Value = 1607
    $ws.Range("F3").Value = 213
    $ws.Range("F4").Value = 201
    $ws.Range("F5").Value = 6098
    $ws.Range("F6").Value = 349
    $ws.Range("F8").Value = 47
    $ws.Range("F10").Value = 8826
    $ws.Range("F11").Value = 2357
    $ws.Range("F12").Value = 256
    $ws.Range("F13").Value = 5415
    $ws.Range("F14").Value = 10300
}
